# Insert two new weekly price rows at the top of the "Acelga" data block
# (rows 880-881), pushing the rest of the dataset down by two rows.
#
# The new rows reuse the Calidad/Volumen/Unidad/Origen/Clasificacion of the
# rows that used to occupy 880-881 (now shifted to 882-883), but carry a new
# reporting date and new price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 880 onward down by 2 rows.
$ws.Rows("880:881").Insert()

# Seed the two new rows from the rows that were pushed down (882 / 883),
# so every column (Mercado, Region, Categoria, Calidad, Volumen, Unidad,
# Origen, Clasificacion, ...) starts out identical to its former self.
$ws.Range("A882:R882").Copy()
$ws.Range("A880:R880").PasteSpecial()

$ws.Range("A883:R883").Copy()
$ws.Range("A881:R881").PasteSpecial()

# Row 880 ("Primera", Volumen 70): new date + new price figures.
$ws.Range("D880").Value = 45218
$ws.Range("K880").Value = 15000
$ws.Range("L880").Value = 15000
$ws.Range("M880").Value = 15000
$ws.Range("P880").Value = 5000

# Row 881 ("Segunda", Volumen 52): new date + new price figures.
$ws.Range("D881").Value = 45218
$ws.Range("K881").Value = 12000
$ws.Range("L881").Value = 12000
$ws.Range("M881").Value = 12000
$ws.Range("P881").Value = 4000
